$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename some test scenario names (column A) to match the updated
# pytest Selenium test cases.
$ws.Range("A5").Value = "Viewing list of posts by category"
$ws.Range("A7").Value = "Leaving a comment with the author more than 60 characters"
$ws.Range("A8").Value = "Leaving a comment without an author"
$ws.Range("A9").Value = "Leaving a comment without a body"

# Reflect the cell that was selected/active when the author last saved
# the workbook.
$ws.Range("A9").Select() | Out-Null
